# Update the symmetric/antisymmetric comparison matrices on the
# "P_valores" and "Estadisticos_DM" sheets with corrected Diebold-Mariano
# values (see commit: "Correcion a Diebold Mariano y revision de Cap1").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: P_valores
# ---------------------------------------------------------------
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.8050953615992988
$wsP.Range("D2").Value = 0.7729222761373844
$wsP.Range("E2").Value = 0.9470163210893376
$wsP.Range("F2").Value = 0.1963803731242795

$wsP.Range("B3").Value = 0.8050953615992988
$wsP.Range("D3").Value = 0.9357673702185372
$wsP.Range("E3").Value = 0.8010203735937627
$wsP.Range("F3").Value = 0.2220054667456803

$wsP.Range("B4").Value = 0.7729222761373844
$wsP.Range("C4").Value = 0.9357673702185372
$wsP.Range("E4").Value = 0.8297782885703899
$wsP.Range("F4").Value = 0.5067971554933928

$wsP.Range("B5").Value = 0.9470163210893376
$wsP.Range("C5").Value = 0.8010203735937627
$wsP.Range("D5").Value = 0.8297782885703899
$wsP.Range("F5").Value = 0.2406585595712734

$wsP.Range("B6").Value = 0.1963803731242795
$wsP.Range("C6").Value = 0.2220054667456803
$wsP.Range("D6").Value = 0.5067971554933928
$wsP.Range("E6").Value = 0.2406585595712734

# ---------------------------------------------------------------
# Sheet: Estadisticos_DM
# ---------------------------------------------------------------
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = -0.2514819579053896
$wsE.Range("D2").Value = -0.2941995757423908
$wsE.Range("E2").Value = -0.06765588580180595
$wsE.Range("F2").Value = -1.356605561128354

$wsE.Range("B3").Value = 0.2514819579053896
$wsE.Range("D3").Value = -0.08205142671232837
$wsE.Range("E3").Value = 0.2568638068843997
$wsE.Range("F3").Value = -1.278081458713226

$wsE.Range("B4").Value = 0.2941995757423908
$wsE.Range("C4").Value = 0.08205142671232837
$wsE.Range("E4").Value = 0.2190423139210384
$wsE.Range("F4").Value = -0.6812974570504345

$wsE.Range("B5").Value = 0.06765588580180595
$wsE.Range("C5").Value = -0.2568638068843997
$wsE.Range("D5").Value = -0.2190423139210384
$wsE.Range("F5").Value = -1.225351358922314

$wsE.Range("B6").Value = 1.356605561128354
$wsE.Range("C6").Value = 1.278081458713226
$wsE.Range("D6").Value = 0.6812974570504345
$wsE.Range("E6").Value = 1.225351358922314
